# Applies the "Updated cryptos list" data refresh to Sheet1.
# All Coin/Link/Price/Volume cells are plain text in this sheet (t="inlineStr"),
# so every assignment below is prefixed with a leading apostrophe to force
# Excel to store it as text instead of auto-coercing to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.946.42"
$ws.Range("E2").Value = "'  -1.13%  "
$ws.Range("D3").Value = "'2.341.33"
$ws.Range("E3").Value = "'  +0.56%  "
$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("D5").Value = "'303.59"
$ws.Range("E5").Value = "'  +0.13%  "
$ws.Range("D6").Value = "'94.49"
$ws.Range("E6").Value = "'  -3.82%  "
$ws.Range("D7").Value = "'0.502"
$ws.Range("E7").Value = "'  -1.08%  "
$ws.Range("E8").Value = "'  +0.05%  "
$ws.Range("D9").Value = "'0.495"
$ws.Range("E9").Value = "'  -1.67%  "
$ws.Range("D10").Value = "'34.13"
$ws.Range("E10").Value = "'  -4.44%  "
$ws.Range("D12").Value = "'18.72"
$ws.Range("E12").Value = "'  -4.55%  "
$ws.Range("D13").Value = "'0.122"
$ws.Range("E13").Value = "'  +2.02%  "
$ws.Range("D14").Value = "'6.71"
$ws.Range("E14").Value = "'  -3.23%  "
$ws.Range("D15").Value = "'2.696.14"
$ws.Range("E15").Value = "'  +0.25%  "
$ws.Range("D16").Value = "'2.322.85"
$ws.Range("E16").Value = "'  +0.25%  "
$ws.Range("D17").Value = "'0.794"
$ws.Range("E17").Value = "'  +0.32%  "
$ws.Range("D18").Value = "'42.898.84"
$ws.Range("E18").Value = "'  -0.93%  "
$ws.Range("D19").Value = "'12.05"
$ws.Range("E19").Value = "'  -5.30%  "
$ws.Range("D20").Value = "'6.22"
$ws.Range("E20").Value = "'  +2.23%  "
$ws.Range("D21").Value = "'0.0₃0888"
$ws.Range("E21").Value = "'  -1.54%  "
$ws.Range("D22").Value = "'67.96"
$ws.Range("D23").Value = "'235.95"
$ws.Range("E23").Value = "'  -0.68%  "
$ws.Range("D24").Value = "'2.21"
$ws.Range("E24").Value = "'  -1.50%  "
$ws.Range("E25").Value = "'  +0.15%  "
$ws.Range("E26").Value = "'  -1.22%  "
$ws.Range("D27").Value = "'24.59"
$ws.Range("E27").Value = "'  -2.00%  "
$ws.Range("E28").Value = "'  -0.23%  "
$ws.Range("E29").Value = "'  -0.24%  "
$ws.Range("D30").Value = "'31.28"
$ws.Range("E30").Value = "'  -6.59%  "
$ws.Range("E31").Value = "'  +0.10%  "
$ws.Range("D32").Value = "'0.0759"
$ws.Range("E32").Value = "'  +7.79%  "
$ws.Range("E33").Value = "'  -0.88%  "
$ws.Range("D34").Value = "'17.24"
$ws.Range("E34").Value = "'  -3.73%  "
$ws.Range("D35").Value = "'4.38"
$ws.Range("E35").Value = "'  -2.92%  "
$ws.Range("E36").Value = "'  -1.08%  "
$ws.Range("B37").Value = "'Monero"
$ws.Range("C37").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "'126.09"
$ws.Range("E37").Value = "'  -23.35%  "
$ws.Range("B38").Value = "'ARBITRUM"
$ws.Range("C38").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "'1.81"
$ws.Range("E38").Value = "'  +2.19%  "
$ws.Range("E39").Value = "'  -0.42%  "
$ws.Range("D40").Value = "'2.76"
$ws.Range("E40").Value = "'  -1.30%  "
$ws.Range("D41").Value = "'22.04"
$ws.Range("E41").Value = "'  +21.16%  "
$ws.Range("E42").Value = "'  -1.34%  "
$ws.Range("D43").Value = "'1.933.08"
$ws.Range("E43").Value = "'  -2.99%  "
$ws.Range("D44").Value = "'0.0281"
$ws.Range("E44").Value = "'  -0.23%  "
$ws.Range("D45").Value = "'10.14"
$ws.Range("E45").Value = "'  -5.55%  "
$ws.Range("E46").Value = "'  +1.11%  "
$ws.Range("E47").Value = "'  -2.94%  "
$ws.Range("B48").Value = "'RocketPoolETH"
$ws.Range("C48").Value = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "'2.566.67"
$ws.Range("E48").Value = "'  +0.41%  "
$ws.Range("B49").Value = "'HuobiToken"
$ws.Range("C49").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").Value = "'2.87"
$ws.Range("E49").Value = "'  -0.77%  "
$ws.Range("D50").Value = "'52.55"
$ws.Range("E50").Value = "'  -2.79%  "
$ws.Range("E51").Value = "'  -1.57%  "
